$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "המאמר היומי של מייק - 21.03.25" "המאמר היומי של מייק - 20.03.25"

Replace-Text "LLMs learn governing principles of dynamical systems, revealing an in-context neural scaling law" "softmax is not enough (for sharp out-of-distribution)"

Replace-Text "המאמר משך את תשומת ליבי כי מופיעים בשמו מודלי שפה ומערכות דינמיות שאני מחבב מהזמנים העליזים של ממבה (state-space models). המאמר טוען שמודלי שפה מפגינים ביצועים טובים בהבנת מערכות דינמיות ממגוון סוגים כולל מערכות סטוכסטיות, כאוטיות, רציפות וכדומה. וכל זה קורה ללא שום טיוב (fine-tune) - כלומר קצת הנדסת פרומפטים ומודל השפה שלכם מבין במערכות דינמיות." "המאמר הזה מציעה שיטה לשיפור ביצועי ההכללה עבור מודלי טרנספורמרים מזווית די לא צפויה. המחברים מציעים שיטה להתמודדות עם מה שנקרא דיספרסיה (או פיזור בעברית) של מקדמים ה-attention בטרנספורמרים. זה מתבטא למשל באי יכולת (לפי המאמר) של הטרנספורמרים למקד את מקדמי ה-attention במספר טוקנים קטן (יחסית לאורך הסדרה). זה חשוב למשל בשאלות כמו מציאת מקסימומים של סדרת מספרים נתונה או שאלות בסגנון ״מחט בערימת השחת״ (needle in a haystack) כאשר המודל מתבקש מקטע קצר לא קשור בטקסט מסוים (יחסית ארוך)."

Replace-Text "לאחרונה היו לא מעט מאמרים שניסו לפצח ״מערכות דינמיות״ הניתנות על ידי דגימות שלהם עם LLMs דרך יצירת דגימות חדשות ממנו (מהמערכת הדינמית) באמצעות LLMs. ההיגיון כאן די פשוט - אם מודל שפה יודע לגנרט מהתפלגות המושרית על ידי מערכת דינמית, אז כנראה הוא מבין אותה. " "המחברים טוענים שאחת הסיבות לבעיות אלו היא פיזור מקדמי ה-attention במנגנון הטרנספורמרים. מקדמים אלו מחושבים עם פונקצית סופטמקס ה״מנרמלת״ את המכפלות הפנימיות של וקטורי K ו-Q עבור כל טוקני הסדרה. לפי המאמר הבעיה קשורה לכך שעבור קונטקסטים ארוכים לסופטמקס במיוחד בטרנספורמרים העמוקים יש ״נטיה למרוח את פלט הסופטמקס״. "

Replace-Text "המחברים לקחו גישה אחרת יותר ישרה - הם הראה שניתן ממש ליצור התפלגות של מערכת דינמית באמצעות LLM כאשר המערכת היא מרקובית. כלומר אם התפלגות דגימה הבאה בזמן t+1 תלויה רק במצב המערכת בזמן t ולא בעבר. עבור מערכת דיסקרטית התפלגות זאת נתונה על ידי מטריצה של הסתברויות מותנות המכילה את ההסתברויות של מצב { x_{t+1} בהינתן מצב x_t בזמן עבור כל הערכים האפשריים שלהם. עבור מערכות רציפות ניתן לבנות מטריצה כזו על ידי דיסקרטיזציה של הערכים של מצבי המערכות." "אחת הדרכים להתמודד עם התופעה הזו היא להוריד את הטמפרטורה אבל זה עלול להעלות סיכוי לשגיאה במקרים בהם הלוגיט (משקל attention לא מנורמל) של הטוקן הנכון יותר קטן מהלוגית המקסימלי. כדי להתמודד עם התופעה המבחרים הציעו גרסה חדשה של סופטמקס בה הטמפרטורה תלויה באנטרופיה של הטוקנים. "

Replace-Text "המאמר מראה שמודלי שפה מצליחים לבנות את מטריצות מעברים בצורה לא רעה במיוחד במצבים שיש יחסית מעט מצבים אפשריים. המרחק בין ההתפלגות החזויה על ידי מודל שפה לבין התפלגות ground truth נמדדה במאמר עם מרחק Bhattacharyya שנתקלתי בו רק בפעם השניה במאמרי deep learning. אציין שהמאמר מציג תוצאות טובות גם עבור מרחקים (divergence) אחרים כמו JSD ו- KL. המחברים מציעים דרך טריקית לבנות את המטריצה הזו עם LLM - מי שרוצה לצלול לעומק, תראו פרק שנקרא Hierarchy-PDF algorithm." "הם אימנו מודל עבור מקרים שבהם הלוגיט של הטוקן הנכון אינו מקסימלי כאשר המטרה היתה למקסם את הסתברות הדגימה של הטוקן הנכון (אחרי מנגנון ה-attention ו-FFN). מטרת המודל היתה לחשב ערך אופטימלי של טמפרטורה כפונקצייה של אנטרופיית של משקלי attention לא מנורמלים. הנוסחה של הטמפרטורה יצאה הופכית (1 חלקי) של פולימום מחזקה 4. אציין כי הטמפרטורה מחושבת בזמן האינפרנס כתלות באינטרופיית הטוקנים לפי המודל הזה."

Replace-Text "וזהו זה, היום זה היה קצר…" "המחברים הראו אמפירית כי עם הטמרטורה האדפטיבית מקטינה פיזור משקלי ה-attention. למרות שהטמפרטורה האדפטיבית האופטימלית יורדת עם עלייה באנטרופיית הלוגיטים היא גורמת לפחות שגיאות של המודל יחסית למקרה שהיא נקבעת באופן קשיח. "

Replace-Text "https://arxiv.org/abs/2402.00795" "https://arxiv.org/abs/2410.01104"
